$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 343.16666
$ws.Range("I11").Value = 343.16666
$ws.Range("K11").Value = 343.16666
$ws.Range("M11").Value = -203.16666
$ws.Range("H76").Value = 4119
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 4158.6665
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 4158.6665
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -4788.6665
$ws.Range("H79").Value = 4119
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 4158.6665
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 4158.6665
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -6342.6665
$ws.Range("H80").Value = 1435.6765
$ws.Range("J80").Value = 1426.3158
$ws.Range("L80").Value = 4278.9474
$ws.Range("N80").Value = -6274.9474
$ws.Range("H83").Value = 1435.6765
$ws.Range("J83").Value = 1426.3158
$ws.Range("L83").Value = 12836.8422
$ws.Range("N83").Value = -22820.8422
$ws.Range("H116").Value = 4580
$ws.Range("H132").Value = 5126.618
$ws.Range("I132").Value = 3200
$ws.Range("J132").Value = 10768.857
$ws.Range("K132").Value = 9600
$ws.Range("L132").Value = 32306.571
$ws.Range("M132").Value = -7070
$ws.Range("N132").Value = -37366.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3850.1233
$ws.Range("I32").Value = 2539.806
$ws.Range("K32").Value = 2539.806
$ws.Range("M32").Value = -2252.806
$ws.Range("H45").Value = 6152.65
$ws.Range("I45").Value = 6499.5356
$ws.Range("K45").Value = 6499.5356
$ws.Range("M45").Value = -6122.5356
$ws.Range("H122").Value = 9999
$ws.Range("I122").Value = 9998.666999999999
$ws.Range("K122").Value = 29996.001
$ws.Range("M122").Value = -27546.001
$ws.Range("H132").Value = 2923.9443
$ws.Range("I132").Value = 2786.4
$ws.Range("J132").Value = 3611.6667
$ws.Range("K132").Value = 8359.200000000001
$ws.Range("L132").Value = 10835.0001
$ws.Range("M132").Value = -5829.200000000001
$ws.Range("N132").Value = -15895.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1755057
$ws.Range("I80").Value = 861
$ws.Range("J80").Value = 2381555.5
$ws.Range("K80").Value = 861
$ws.Range("L80").Value = 2381555.5
$ws.Range("M80").Value = 137
$ws.Range("N80").Value = -2383551.5
$ws.Range("H83").Value = 1755057
$ws.Range("I83").Value = 861
$ws.Range("J83").Value = 2381555.5
$ws.Range("K83").Value = 4305
$ws.Range("L83").Value = 11907777.5
$ws.Range("M83").Value = 687
$ws.Range("N83").Value = -11917761.5
$ws.Range("H86").Value = 2070
$ws.Range("I86").Value = 2068.3845
$ws.Range("K86").Value = 2068.3845
$ws.Range("M86").Value = -945.3845000000001
$ws.Range("H89").Value = 2070
$ws.Range("I89").Value = 2068.3845
$ws.Range("K89").Value = 10341.9225
$ws.Range("M89").Value = -4725.922500000001
$ws.Range("H95").Value = 24966
$ws.Range("J95").Value = 24966
$ws.Range("L95").Value = 24966
$ws.Range("N95").Value = -30458

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2604.4
$ws.Range("I31").Value = 2604.4
$ws.Range("K31").Value = 2604.4
$ws.Range("M31").Value = -2309.4
$ws.Range("H34").Value = 2604.4
$ws.Range("I34").Value = 2604.4
$ws.Range("K34").Value = 2604.4
$ws.Range("M34").Value = -2402.4
$ws.Range("H41").Value = 17000
$ws.Range("H51").Value = 37385.625
$ws.Range("I51").Value = 17818
$ws.Range("J51").Value = 69998.336
$ws.Range("K51").Value = 17818
$ws.Range("L51").Value = 69998.336
$ws.Range("M51").Value = -17082
$ws.Range("N51").Value = -71470.336
$ws.Range("H60").Value = 52562.363
$ws.Range("I60").Value = 27273.25
$ws.Range("J60").Value = 120000
$ws.Range("K60").Value = 27273.25
$ws.Range("L60").Value = 120000
$ws.Range("M60").Value = -26762.25
$ws.Range("N60").Value = -121022
$ws.Range("H61").Value = 37385.625
$ws.Range("I61").Value = 17818
$ws.Range("J61").Value = 69998.336
$ws.Range("K61").Value = 17818
$ws.Range("L61").Value = 69998.336
$ws.Range("M61").Value = -17470
$ws.Range("N61").Value = -70694.336
$ws.Range("H74").Value = 44812.445
$ws.Range("J74").Value = 44812.445
$ws.Range("L74").Value = 44812.445
$ws.Range("N74").Value = -46560.445
$ws.Range("H77").Value = 44812.445
$ws.Range("J77").Value = 44812.445
$ws.Range("L77").Value = 134437.335
$ws.Range("N77").Value = -143173.335
$ws.Range("H86").Value = 4766447
$ws.Range("I86").Value = 6670755.5
$ws.Range("K86").Value = 6670755.5
$ws.Range("M86").Value = -6669632.5
$ws.Range("H88").Value = 9177.5
$ws.Range("J88").Value = 9177.5
$ws.Range("L88").Value = 9177.5
$ws.Range("N88").Value = -9989.5
$ws.Range("H89").Value = 4766447
$ws.Range("I89").Value = 6670755.5
$ws.Range("K89").Value = 33353777.5
$ws.Range("M89").Value = -33348161.5
$ws.Range("H91").Value = 9177.5
$ws.Range("J91").Value = 9177.5
$ws.Range("L91").Value = 9177.5
$ws.Range("N91").Value = -11985.5
$ws.Range("H94").Value = 5086.125
$ws.Range("I94").Value = 6532.778
$ws.Range("J94").Value = 3226.1428
$ws.Range("K94").Value = 6532.778
$ws.Range("L94").Value = 3226.1428
$ws.Range("M94").Value = -6081.778
$ws.Range("N94").Value = -4128.1428
$ws.Range("H105").Value = 1727.9333
$ws.Range("I105").Value = 1811.6666
$ws.Range("J105").Value = 1393
$ws.Range("K105").Value = 1811.6666
$ws.Range("L105").Value = 1393
$ws.Range("M105").Value = -64.66660000000002
$ws.Range("N105").Value = -4887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 1013.1667
$ws.Range("I24").Value = 324
$ws.Range("K24").Value = 972
$ws.Range("M24").Value = -742
$ws.Range("H45").Value = 3069.4
$ws.Range("I45").Value = 6365
$ws.Range("J45").Value = 872.3333
$ws.Range("K45").Value = 19095
$ws.Range("L45").Value = 2616.9999
$ws.Range("M45").Value = -18563
$ws.Range("N45").Value = -3680.9999
$ws.Range("H50").Value = 1341.3636
$ws.Range("I50").Value = 1466.6666
$ws.Range("J50").Value = 777.5
$ws.Range("K50").Value = 4399.9998
$ws.Range("L50").Value = 2332.5
$ws.Range("M50").Value = -3918.9998
$ws.Range("N50").Value = -3294.5
$ws.Range("H53").Value = 1341.3636
$ws.Range("I53").Value = 1466.6666
$ws.Range("J53").Value = 777.5
$ws.Range("K53").Value = 4399.9998
$ws.Range("L53").Value = 2332.5
$ws.Range("M53").Value = -3918.9998
$ws.Range("N53").Value = -3294.5
$ws.Range("H63").Value = 12111
$ws.Range("I63").Value = 7777.5
$ws.Range("K63").Value = 23332.5
$ws.Range("M63").Value = -22583.5
$ws.Range("H66").Value = 12111
$ws.Range("I66").Value = 7777.5
$ws.Range("K66").Value = 69997.5
$ws.Range("M66").Value = -66253.5
$ws.Range("H131").Value = 7048349.5
$ws.Range("I131").Value = 11029901
$ws.Range("J131").Value = 5057574
$ws.Range("K131").Value = 33089703
$ws.Range("L131").Value = 15172722
$ws.Range("M131").Value = -33084663
$ws.Range("N131").Value = -15182802

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4428.3335
$ws.Range("I102").Value = 4731.8823
$ws.Range("J102").Value = 3138.25
$ws.Range("K102").Value = 4731.8823
$ws.Range("L102").Value = 3138.25
$ws.Range("M102").Value = -3109.8823
$ws.Range("N102").Value = -6382.25
$ws.Range("H122").Value = 559300.25
$ws.Range("I122").Value = 1113847.2
$ws.Range("J122").Value = 4753.222
$ws.Range("K122").Value = 3341541.6
$ws.Range("L122").Value = 14259.666
$ws.Range("M122").Value = -3339091.6
$ws.Range("N122").Value = -19159.666
$ws.Range("H126").Value = 6786.4707
$ws.Range("J126").Value = 6898.3335
$ws.Range("L126").Value = 20695.0005
$ws.Range("N126").Value = -25635.0005
$ws.Range("H132").Value = 4267.5
$ws.Range("I132").Value = 1951.5
$ws.Range("J132").Value = 7162.5
$ws.Range("K132").Value = 5854.5
$ws.Range("L132").Value = 21487.5
$ws.Range("M132").Value = -3324.5
$ws.Range("N132").Value = -26547.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3568.7646
$ws.Range("I40").Value = 2579.6365
$ws.Range("J40").Value = 5382.1665
$ws.Range("K40").Value = 2579.6365
$ws.Range("L40").Value = 5382.1665
$ws.Range("M40").Value = -2443.6365
$ws.Range("N40").Value = -5654.1665
$ws.Range("H106").Value = 20331.666
$ws.Range("J106").Value = 20331.666
$ws.Range("L106").Value = 20331.666
$ws.Range("N106").Value = -22855.666
$ws.Range("H122").Value = 7115.5835
$ws.Range("I122").Value = 5527.4287
$ws.Range("K122").Value = 16582.2861
$ws.Range("M122").Value = -14132.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7072.909
$ws.Range("I122").Value = 4649.68
$ws.Range("J122").Value = 14645.5
$ws.Range("K122").Value = 13949.04
$ws.Range("L122").Value = 43936.5
$ws.Range("M122").Value = -11499.04
$ws.Range("N122").Value = -48836.5
$ws.Range("H136").Value = 8284.317999999999
$ws.Range("I136").Value = 10114.615
$ws.Range("J136").Value = 5640.5557
$ws.Range("K136").Value = 30343.845
$ws.Range("L136").Value = 16921.6671
$ws.Range("M136").Value = -27793.845
$ws.Range("N136").Value = -22021.6671
